{"js": "// Office.js (Word JavaScript API) script\n// Applies the edit described by the diff:\n//  1. TITLE paragraph: new subtitle text (\"Abstractive Summarization of\n//     Threat Intelligence using Transformers\") replacing the old one,\n//     keeping \"TITLE: \" bold and the rest not bold.\n//  2. Proposer paragraph: merge \"Proposer:\" + \" \" into a single\n//     \"Proposer: \" run (no visible text change).\n//  3. DELIVERABLE paragraph: replace the deliverable description.\n//  4. SKILLS NEEDED paragraph: replace the tool list tail.\n\nconst body = context.document.body;\n\n// --- 1. TITLE ---------------------------------------------------------\n// Merge \"TITLE\" + \": \" into a single bold run reading \"TITLE: \".\nlet results = body.search(\"TITLE: \", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"TITLE: \", \"Replace\");\n}\nawait context.sync();\n\n// Replace the old (non-bold) subtitle text with the new project title.\nresults = body.search(\n  \"Classification of threat intelligence from news articles using NLP\",\n  { matchCase: true }\n);\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\n    \"Abstractive Summarization of Threat Intelligence using Transformers\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// --- 2. Proposer --------------------------------------------------------\n// Merge \"Proposer:\" + \" \" into a single bold run \"Proposer: \".\nresults = body.search(\"Proposer: \", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"Proposer: \", \"Replace\");\n}\nawait context.sync();\n\n// --- 3. DELIVERABLE ------------------------------------------------------\nresults = body.search(\n  \"An ML pipeline to perform classification task on NLP dataset \",\n  { matchCase: true }\n);\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\n    \"A utility which takes input (text) and generates its summary and relevant keywords by using Transformer based deep learning models.\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// --- 4. SKILLS NEEDED ------------------------------------------------------\nresults = body.search(\", Python3, Keras/TensorFlow, Jupyter Notebook\", {\n  matchCase: true,\n});\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\", Python3, Torch, Flask\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Applies the edit described by the diff:\n#   1. TITLE paragraph: new subtitle text (\"Abstractive Summarization of\n#      Threat Intelligence using Transformers\") replacing the old one,\n#      keeping \"TITLE: \" bold and the rest not bold.\n#   2. Proposer paragraph: merge \"Proposer:\" + \" \" into a single\n#      \"Proposer: \" run (no visible text change).\n#   3. DELIVERABLE paragraph: replace the deliverable description.\n#   4. SKILLS NEEDED paragraph: replace the tool list tail.\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText($oldText, $newText) {\n    $f = $d.Content.Find\n    $f.Text = $oldText\n    $f.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\n# --- 1. TITLE -------------------------------------------------------------\n# Merge \"TITLE\" + \": \" into a single bold run reading \"TITLE: \".\nReplace-DocText \"TITLE: \" \"TITLE: \"\n# Replace the old (non-bold) subtitle text with the new project title.\nReplace-DocText \"Classification of threat intelligence from news articles using NLP\" \"Abstractive Summarization of Threat Intelligence using Transformers\"\n\n# --- 2. Proposer ------------------------------------------------------------\n# Merge \"Proposer:\" + \" \" into a single bold run \"Proposer: \".\nReplace-DocText \"Proposer: \" \"Proposer: \"\n\n# --- 3. DELIVERABLE ---------------------------------------------------------\nReplace-DocText \"An ML pipeline to perform classification task on NLP dataset \" \"A utility which takes input (text) and generates its summary and relevant keywords by using Transformer based deep learning models.\"\n\n# --- 4. SKILLS NEEDED --------------------------------------------------------\nReplace-DocText \", Python3, Keras/TensorFlow, Jupyter Notebook\" \", Python3, Torch, Flask\"\n"}
